$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending / target cluster labels used in this table (cross product of 3 clusters)
$clusters = @("ECs", "FAPs", "sCs")
$ligand = "Efnb3"
$receptor = "Ephb4"

# Numeric columns E..T for each (sending, target) cluster combination.
# Row order mirrors the cross-product: ECs x {ECs,FAPs,sCs}, FAPs x {...}, sCs x {...}
$data = @(
    @(1, 0.3333333333333333, 0.071327, 0.213981, 0.03356605248408491, 0.03356605248408491, 3, 1, 17.00372766666667, 51.011183, 0.6620200065567142, 0.6620200065567141, 1.212824883280333, 10.915423949523, 0.02222139828559691, 0.0222213982855969),
    @(1, 0.3333333333333333, 0.071327, 0.213981, 0.03356605248408491, 0.03356605248408491, 3, 1, 5.393811333333333, 16.181434, 0.2100016587103467, 0.2100016587103467, 0.3847243809726666, 3.462519428754, 0.007048926698016383, 0.007048926698016384),
    @(1, 0.3333333333333333, 0.071327, 0.213981, 0.03356605248408491, 0.03356605248408491, 3, 1, 3.287074, 9.861222, 0.1279783347329392, 0.1279783347329391, 0.234457127198, 2.110114144782, 0.004295727500471623, 0.004295727500471622),
    @(3, 1, 0.308223, 0.9246690000000001, 0.1450478695977975, 0.1450478695977975, 3, 1, 17.00372766666667, 51.011183, 0.6620200065567142, 0.6620200065567141, 5.240939952603, 47.168459573427, 0.09602459158217133, 0.09602459158217132),
    @(3, 1, 0.308223, 0.9246690000000001, 0.1450478695977975, 0.1450478695977975, 3, 1, 5.393811333333333, 16.181434, 0.2100016587103467, 0.2100016587103467, 1.662496710594, 14.962470395346, 0.03046029320793954, 0.03046029320793954),
    @(3, 1, 0.308223, 0.9246690000000001, 0.1450478695977975, 0.1450478695977975, 3, 1, 3.287074, 9.861222, 0.1279783347329392, 0.1279783347329391, 1.013151809502, 9.118366285518, 0.01856298480768664, 0.01856298480768663),
    @(3, 1, 1.745424333333333, 5.236273, 0.8213860779181176, 0.8213860779181176, 3, 1, 17.00372766666667, 51.011183, 0.6620200065567142, 0.6620200065567141, 29.67872002677322, 267.108480240959, 0.543774016688946, 0.5437740166889459),
    @(3, 1, 1.745424333333333, 5.236273, 0.8213860779181176, 0.8213860779181176, 3, 1, 5.393811333333333, 16.181434, 0.2100016587103467, 0.2100016587103467, 9.41448955060911, 84.730405955482, 0.1724924388043907, 0.1724924388043907),
    @(3, 1, 1.745424333333333, 5.236273, 0.8213860779181176, 0.8213860779181176, 3, 1, 3.287074, 9.861222, 0.1279783347329392, 0.1279783347329391, 5.737338945067333, 51.63605050560599, 0.1051196224247809, 0.1051196224247809)
)

$row = 2
for ($i = 0; $i -lt 3; $i++) {
    for ($j = 0; $j -lt 3; $j++) {
        $ws.Cells.Item($row, 1).Value = $clusters[$i]
        $ws.Cells.Item($row, 2).Value = $ligand
        $ws.Cells.Item($row, 3).Value = $receptor
        $ws.Cells.Item($row, 4).Value = $clusters[$j]

        $vals = $data[($i * 3) + $j]
        for ($c = 0; $c -lt $vals.Length; $c++) {
            $ws.Cells.Item($row, 5 + $c).Value = $vals[$c]
        }

        $row++
    }
}
